$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new daily rows (2025-11-10) for both stations, continuing the
# existing log that currently ends at row 19.
$ws.Range("A20").Value = 45971
$ws.Range("B20").Value = "四方坪站"
$ws.Range("C20").Value = 10798.15
$ws.Range("D20").Value = 9521.94
$ws.Range("E20").Value = 3568.03
$ws.Range("F20").Value = 441

$ws.Range("A21").Value = 45971
$ws.Range("B21").Value = "高岭站"
$ws.Range("C21").Value = 4554.88
$ws.Range("D21").Value = 4173.13
$ws.Range("E21").Value = 1219.25
$ws.Range("F21").Value = 167

$ws.Range("G20").Select()
